$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Set the new cell values for columns D and E (rows 2-6) ---
# Values are entered in the same order the original author typed them in,
# so the shared-string table ends up built in the same sequence.
$ws.Range("D2").Value = "数据类型"
$ws.Range("D4").Value = "type"
$ws.Range("E2").Value = "数据class"
$ws.Range("E4").Value = "clz"
$ws.Range("E5").Value = "Integer"
$ws.Range("D5").Value = "OBJECT"
$ws.Range("D6").Value = "ARRAY"
$ws.Range("E6").Value = "Integer"

$ws.Range("D3").Value = "string"
$ws.Range("E3").Value = "string"

# --- Copy the header/row formatting from column C into the new D/E columns ---
# Rows 2-4 carry explicit cell styles in columns A-C (header/sub-header/key
# row styles); rows 5-6 use the default (unstyled) format already, so only
# rows 2, 3 & 4 need their formatting copied across.
$ws.Range("C2").Copy()
$ws.Range("D2:E2").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C3").Copy()
$ws.Range("D3:E3").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C4").Copy()
$ws.Range("D4:E4").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# --- Update the active selection on the sheet to match the new edit location ---
$ws.Range("D6").Select()
